# Auto-generated edit script applying value updates described by the commit diff
# to Sheets/Raiden_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1250057.6
$ws.Range("J9").Value = 2500074.5
$ws.Range("L9").Value = 2500074.5
$ws.Range("N9").Value = -2500412.5

$ws.Range("H40").Value = 4374
$ws.Range("I40").Value = 4361.1665
$ws.Range("K40").Value = 4361.1665
$ws.Range("M40").Value = -4186.1665

$ws.Range("H42").Value = 481.125
$ws.Range("J42").Value = 597.8
$ws.Range("L42").Value = 1793.4
$ws.Range("N42").Value = -2253.4

$ws.Range("H98").Value = 2045.9445
$ws.Range("I98").Value = 1712.4546
$ws.Range("K98").Value = 1712.4546
$ws.Range("M98").Value = -214.4546

$ws.Range("H112").Value = 1844.5186
$ws.Range("J112").Value = 1864.44
$ws.Range("L112").Value = 5593.32
$ws.Range("N112").Value = -7809.32

$ws.Range("H122").Value = 2045.9445
$ws.Range("I122").Value = 1712.4546
$ws.Range("K122").Value = 5137.3638
$ws.Range("M122").Value = -2687.3638

$ws.Range("H132").Value = 240423.36
$ws.Range("I132").Value = 2186.7144
$ws.Range("K132").Value = 6560.1432
$ws.Range("M132").Value = -4030.1432

$ws.Range("H137").Value = 3838.6428
$ws.Range("I137").Value = 2701.4285
$ws.Range("J137").Value = 4975.857
$ws.Range("K137").Value = 8104.2855
$ws.Range("L137").Value = 14927.571
$ws.Range("M137").Value = -5554.2855
$ws.Range("N137").Value = -20027.571

$ws.Range("H138").Value = 2068.4707
$ws.Range("I138").Value = 2381.5
$ws.Range("J138").Value = 1938.0416
$ws.Range("K138").Value = 7144.5
$ws.Range("L138").Value = 5814.1248
$ws.Range("M138").Value = -2004.5
$ws.Range("N138").Value = -16094.1248

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2536.2222
$ws.Range("I45").Value = 2379.8333
$ws.Range("K45").Value = 2379.8333
$ws.Range("M45").Value = -2002.8333

$ws.Range("H61").Value = 3690.36
$ws.Range("I61").Value = 3102.353
$ws.Range("K61").Value = 3102.353
$ws.Range("M61").Value = -2890.353

$ws.Range("H75").Value = 75625.75
$ws.Range("J75").Value = 75625.75
$ws.Range("L75").Value = 75625.75
$ws.Range("N75").Value = -77373.75

$ws.Range("H78").Value = 75625.75
$ws.Range("J78").Value = 75625.75
$ws.Range("L78").Value = 226877.25
$ws.Range("N78").Value = -235613.25

$ws.Range("H80").Value = 28551.666
$ws.Range("J80").Value = 37827.5
$ws.Range("L80").Value = 37827.5
$ws.Range("N80").Value = -39823.5

$ws.Range("H83").Value = 28551.666
$ws.Range("J83").Value = 37827.5
$ws.Range("L83").Value = 113482.5
$ws.Range("N83").Value = -123466.5

$ws.Range("H122").Value = 2424
$ws.Range("J122").Value = 3497
$ws.Range("L122").Value = 10491
$ws.Range("N122").Value = -15391

$ws.Range("H136").Value = 3690.36
$ws.Range("I136").Value = 3102.353
$ws.Range("K136").Value = 9307.059000000001
$ws.Range("M136").Value = -6757.059000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2515.6
$ws.Range("J80").Value = 2661.111
$ws.Range("L80").Value = 2661.111
$ws.Range("N80").Value = -4657.111

$ws.Range("H83").Value = 2515.6
$ws.Range("J83").Value = 2661.111
$ws.Range("L83").Value = 13305.555
$ws.Range("N83").Value = -23289.555

$ws.Range("H86").Value = 2769.4443
$ws.Range("I86").Value = 2791
$ws.Range("K86").Value = 2791
$ws.Range("M86").Value = -1668

$ws.Range("H89").Value = 2769.4443
$ws.Range("I89").Value = 2791
$ws.Range("K89").Value = 13955
$ws.Range("M89").Value = -8339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1790
$ws.Range("I16").Value = 1784
$ws.Range("J16").Value = 1796
$ws.Range("K16").Value = 1784
$ws.Range("L16").Value = 1796
$ws.Range("M16").Value = -1497
$ws.Range("N16").Value = -2370

$ws.Range("H58").Value = 1702.975
$ws.Range("I58").Value = 1659.4445
$ws.Range("K58").Value = 1659.4445
$ws.Range("M58").Value = -1456.4445

$ws.Range("H74").Value = 41406.09
$ws.Range("J74").Value = 41406.09
$ws.Range("L74").Value = 41406.09
$ws.Range("N74").Value = -43154.09

$ws.Range("H77").Value = 41406.09
$ws.Range("J77").Value = 41406.09
$ws.Range("L77").Value = 124218.27
$ws.Range("N77").Value = -132954.27

$ws.Range("H88").Value = 42397.75
$ws.Range("J88").Value = 39026
$ws.Range("L88").Value = 39026
$ws.Range("N88").Value = -39838

$ws.Range("H91").Value = 42397.75
$ws.Range("J91").Value = 39026
$ws.Range("L91").Value = 39026
$ws.Range("N91").Value = -41834

$ws.Range("H113").Value = 1790
$ws.Range("I113").Value = 1784
$ws.Range("J113").Value = 1796
$ws.Range("K113").Value = 1784
$ws.Range("L113").Value = 1796
$ws.Range("M113").Value = 386
$ws.Range("N113").Value = -6136

$ws.Range("H132").Value = 3871.0454
$ws.Range("J132").Value = 3685.3333
$ws.Range("L132").Value = 11055.9999
$ws.Range("N132").Value = -16115.9999

$ws.Range("H136").Value = 1702.975
$ws.Range("I136").Value = 1659.4445
$ws.Range("K136").Value = 4978.333500000001
$ws.Range("M136").Value = -2428.333500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 350
$ws.Range("J12").Value = 299
$ws.Range("L12").Value = 897
$ws.Range("N12").Value = -1243

$ws.Range("H51").Value = 2321.6
$ws.Range("J51").Value = 2279.5
$ws.Range("L51").Value = 6838.5
$ws.Range("N51").Value = -7758.5

$ws.Range("H76").Value = 3874.5
$ws.Range("I76").Value = 499
$ws.Range("J76").Value = 7250
$ws.Range("K76").Value = 1497
$ws.Range("L76").Value = 21750
$ws.Range("M76").Value = -1114
$ws.Range("N76").Value = -22516

$ws.Range("H79").Value = 3874.5
$ws.Range("I79").Value = 499
$ws.Range("J79").Value = 7250
$ws.Range("K79").Value = 1497
$ws.Range("L79").Value = 21750
$ws.Range("M79").Value = -171
$ws.Range("N79").Value = -24402

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2203.5667
$ws.Range("I102").Value = 1710.95
$ws.Range("J102").Value = 3188.8
$ws.Range("K102").Value = 1710.95
$ws.Range("L102").Value = 3188.8
$ws.Range("M102").Value = -88.95000000000005
$ws.Range("N102").Value = -6432.8

$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988

$ws.Range("H113").Value = 3281.2354
$ws.Range("I113").Value = 2199
$ws.Range("K113").Value = 2199
$ws.Range("M113").Value = -29

$ws.Range("H122").Value = 47551.684
$ws.Range("I122").Value = 54660.42
$ws.Range("J122").Value = 2529.6667
$ws.Range("K122").Value = 163981.26
$ws.Range("L122").Value = 7589.000100000001
$ws.Range("M122").Value = -161531.26
$ws.Range("N122").Value = -12489.0001

$ws.Range("H132").Value = 3012.739
$ws.Range("I132").Value = 2309.4736
$ws.Range("K132").Value = 6928.4208
$ws.Range("M132").Value = -4398.4208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3709.889
$ws.Range("I7").Value = 2678.6
$ws.Range("K7").Value = 2678.6
$ws.Range("M7").Value = -2566.6

$ws.Range("H40").Value = 7498.1816
$ws.Range("I40").Value = 7783.7144
$ws.Range("K40").Value = 7783.7144
$ws.Range("M40").Value = -7647.7144

$ws.Range("H46").Value = 1057.6666
$ws.Range("I46").Value = 1057.6666
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1057.6666
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -869.6666
$ws.Range("N46").ClearContents()

$ws.Range("H122").Value = 6244.25
$ws.Range("I122").Value = 3980
$ws.Range("J122").Value = 6999
$ws.Range("K122").Value = 11940
$ws.Range("L122").Value = 20997
$ws.Range("M122").Value = -9490
$ws.Range("N122").Value = -25897

$ws.Range("H126").Value = 3709.889
$ws.Range("I126").Value = 2678.6
$ws.Range("K126").Value = 8035.799999999999
$ws.Range("M126").Value = -5565.799999999999

$ws.Range("H132").Value = 3593.2354
$ws.Range("I132").Value = 3432.6667
$ws.Range("K132").Value = 10298.0001
$ws.Range("M132").Value = -7768.000100000001

$ws.Range("H133").Value = 66748
$ws.Range("J133").Value = 66748
$ws.Range("L133").Value = 66748
$ws.Range("N133").Value = -71808

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 15999.75
$ws.Range("I4").Value = 16666.334
$ws.Range("J4").Value = 14000
$ws.Range("K4").Value = 16666.334
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = -16553.334
$ws.Range("N4").Value = -14226

$ws.Range("H122").Value = 3406.739
$ws.Range("I122").Value = 2690.9167
$ws.Range("K122").Value = 8072.750100000001
$ws.Range("M122").Value = -5622.750100000001

$ws.Range("H126").Value = 2565.96
$ws.Range("I126").Value = 2489.1304
$ws.Range("K126").Value = 7467.3912
$ws.Range("M126").Value = -4997.3912

$ws.Range("H132").Value = 4602.12
$ws.Range("I132").Value = 4727.625
$ws.Range("J132").Value = 1590
$ws.Range("K132").Value = 14182.875
$ws.Range("L132").Value = 4770
$ws.Range("M132").Value = -11652.875
$ws.Range("N132").Value = -9830
